$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 6.673343051422484
$ws.Range("D2").Value = 5.734513289306478
$ws.Range("E2").Value = 11.31826408210318
$ws.Range("F2").Value = 42.14967022228711
$ws.Range("G2").Value = 52.80682490552394
$ws.Range("H2").Value = 20.31745731006703
$ws.Range("I2").Value = 32.79796680540149
$ws.Range("J2").Value = 9.529443551007709
$ws.Range("L2").Value = 9.375771395673613
$ws.Range("M2").Value = 61.43846333124413

$ws.Range("C3").Value = 6.681041832947805
$ws.Range("D3").Value = 5.687414038674131
$ws.Range("E3").Value = 11.2613235182166
$ws.Range("F3").Value = 42.74915512258869
$ws.Range("G3").Value = 53.53667225040169
$ws.Range("H3").Value = 20.54675264338139
$ws.Range("I3").Value = 33.23912315327222
$ws.Range("J3").Value = 9.569899224022988
$ws.Range("L3").Value = 9.372842443310841
$ws.Range("M3").Value = 58.23483566748132

$ws.Range("C4").Value = 6.686801126405119
$ws.Range("D4").Value = 5.658205573558097
$ws.Range("E4").Value = 11.22932914586932
$ws.Range("F4").Value = 43.1390554449941
$ws.Range("G4").Value = 54.02072405493974
$ws.Range("H4").Value = 20.6955840357793
$ws.Range("I4").Value = 33.52559033610807
$ws.Range("J4").Value = 9.597555175733657
$ws.Range("L4").Value = 9.373477962316091
$ws.Range("M4").Value = 56.16797407622531

$ws.Range("C5").Value = 6.689407651803537
$ws.Range("D5").Value = 5.646235863479674
$ws.Range("E5").Value = 11.21704091147557
$ws.Range("F5").Value = 43.30336136916682
$ws.Range("G5").Value = 54.22681461200769
$ws.Range("H5").Value = 20.75824091056126
$ws.Range("I5").Value = 33.64620907896523
$ws.Range("J5").Value = 9.609525829153192
$ws.Range("L5").Value = 9.37434494655623
$ws.Range("M5").Value = 55.30108309999704

$ws.Range("C6").Value = 6.689856140926112
$ws.Range("D6").Value = 5.644244422177898
$ws.Range("E6").Value = 11.2150458337274
$ws.Range("F6").Value = 43.33096962518434
$ws.Range("G6").Value = 54.26156384491723
$ws.Range("H6").Value = 20.76876580798456
$ws.Range("I6").Value = 33.66647099991461
$ws.Range("J6").Value = 9.611555650587594
$ws.Range("L6").Value = 9.374525512062732
$ws.Range("M6").Value = 55.1556649816991

$ws.Range("C7").Value = 6.686835227918285
$ws.Range("D7").Value = 5.658044409650068
$ws.Range("E7").Value = 11.2291603825272
$ws.Range("F7").Value = 43.14124947845234
$ws.Range("G7").Value = 54.0234679433152
$ws.Range("H7").Value = 20.69642094245073
$ws.Range("I7").Value = 33.52720137997697
$ws.Range("J7").Value = 9.597713788967779
$ws.Range("L7").Value = 9.3734871983141
$ws.Range("M7").Value = 56.15638189652078

$ws.Range("C8").Value = 6.675783356686321
$ws.Range("D8").Value = 5.718336944361003
$ws.Range("E8").Value = 11.29801740182394
$ws.Range("F8").Value = 42.35180345934765
$ws.Range("G8").Value = 53.05089875970305
$ws.Range("H8").Value = 20.39483930904912
$ws.Range("I8").Value = 32.94681549405306
$ws.Range("J8").Value = 9.542803609186384
$ws.Range("L8").Value = 9.374253923374026
$ws.Range("M8").Value = 60.35482201737097

$ws.Range("C9").Value = 6.662302300687874
$ws.Range("D9").Value = 5.834026725980261
$ws.Range("E9").Value = 11.45641451822917
$ws.Range("F9").Value = 40.98001310785841
$ws.Range("G9").Value = 51.43806692976239
$ws.Range("H9").Value = 19.86800143011201
$ws.Range("I9").Value = 31.93437262008479
$ws.Range("J9").Value = 9.457820782860441
$ws.Range("L9").Value = 9.395251550708334
$ws.Range("M9").Value = 67.78013527388318

$ws.Range("C10").Value = 6.65739655994985
$ws.Range("D10").Value = 5.917147617584348
$ws.Range("E10").Value = 11.58677636608921
$ws.Range("F10").Value = 40.08445656082965
$ws.Range("G10").Value = 50.4461904755676
$ws.Range("H10").Value = 19.52139066273046
$ws.Range("I10").Value = 31.27006268304687
$ws.Range("J10").Value = 9.409729517744063
$ws.Range("L10").Value = 9.42281588868917
$ws.Range("M10").Value = 72.72772105336047

$ws.Range("C11").Value = 6.656251688412422
$ws.Range("D11").Value = 5.954486387083196
$ws.Range("E11").Value = 11.64905505046076
$ws.Range("F11").Value = 39.70260467454466
$ws.Range("G11").Value = 50.04006908227124
$ws.Range("H11").Value = 19.37275925630432
$ws.Range("I11").Value = 30.98583113593287
$ws.Range("J11").Value = 9.39108410197402
$ws.Range("L11").Value = 9.438045014415234
$ws.Range("M11").Value = 74.86633959651576

$ws.Range("C12").Value = 6.655974531167091
$ws.Range("D12").Value = 5.968551917030519
$ws.Range("E12").Value = 11.67306023029652
$ws.Range("F12").Value = 39.56178924431418
$ws.Range("G12").Value = 49.89304734508645
$ws.Range("H12").Value = 19.31780162673316
$ws.Range("I12").Value = 30.88084984459807
$ws.Range("J12").Value = 9.384498508675591
$ws.Range("L12").Value = 9.444202925236626
$ws.Range("M12").Value = 75.65996881028816

$ws.Range("C13").Value = 6.656027265053488
$ws.Range("D13").Value = 5.965526028935384
$ws.Range("E13").Value = 11.66787163941568
$ws.Range("F13").Value = 39.59194623069823
$ws.Range("G13").Value = 49.92440544281484
$ws.Range("H13").Value = 19.32957834254636
$ws.Range("I13").Value = 30.90334040539993
$ws.Range("J13").Value = 9.385895538247722
$ws.Range("L13").Value = 9.442859254439799
$ws.Range("M13").Value = 75.48976914020803

$ws.Range("C14").Value = 6.656225751799242
$ws.Range("D14").Value = 5.955645079523954
$ws.Range("E14").Value = 11.65102155933261
$ws.Range("F14").Value = 39.69094320939008
$ws.Range("G14").Value = 50.02783620579879
$ws.Range("H14").Value = 19.3682111346081
$ws.Range("I14").Value = 30.97714072542884
$ws.Range("J14").Value = 9.390532720483987
$ws.Range("L14").Value = 9.438543770120695
$ws.Range("M14").Value = 74.93195745835089

$ws.Range("C15").Value = 6.656367699341955
$ws.Range("D15").Value = 5.949582925303137
$ws.Range("E15").Value = 11.64075510815438
$ws.Range("F15").Value = 39.75207789261471
$ws.Range("G15").Value = 50.09208054724987
$ws.Range("H15").Value = 19.39204828775973
$ws.Range("I15").Value = 31.02269289072506
$ws.Range("J15").Value = 9.393435304104141
$ws.Range("L15").Value = 9.435951450298887
$ws.Range("M15").Value = 74.58816708375709

$ws.Range("C16").Value = 6.657493257237085
$ws.Range("D16").Value = 5.914697419140227
$ws.Range("E16").Value = 11.58276565849896
$ws.Range("F16").Value = 40.10993470458386
$ws.Range("G16").Value = 50.47366270880151
$ws.Range("H16").Value = 19.53128818267558
$ws.Range("I16").Value = 31.28900505732279
$ws.Range("J16").Value = 9.411013939532856
$ws.Range("L16").Value = 9.42187511447673
$ws.Range("M16").Value = 72.58569238351117

$ws.Range("C17").Value = 6.658462169072365
$ws.Range("D17").Value = 5.893170904504711
$ws.Range("E17").Value = 11.54794886435324
$ws.Range("F17").Value = 40.33608605990852
$ws.Range("G17").Value = 50.71949678136004
$ws.Range("H17").Value = 19.61904110225014
$ws.Range("I17").Value = 31.45702687740077
$ws.Range("J17").Value = 9.422632751557217
$ws.Range("L17").Value = 9.413931669912548
$ws.Range("M17").Value = 71.32844870771694

$ws.Range("C18").Value = 6.659121746905148
$ws.Range("D18").Value = 5.880745208809119
$ws.Range("E18").Value = 11.52820372365369
$ws.Range("F18").Value = 40.46855925300766
$ws.Range("G18").Value = 50.8651200896337
$ws.Range("H18").Value = 19.67036385197242
$ws.Range("I18").Value = 31.55535530597948
$ws.Range("J18").Value = 9.429619108711245
$ws.Range("L18").Value = 9.409615756872505
$ws.Range("M18").Value = 70.59477076896006

$ws.Range("C19").Value = 6.659362633296824
$ws.Range("D19").Value = 5.876530666144155
$ws.Range("E19").Value = 11.52156673560272
$ws.Range("F19").Value = 40.51382084009072
$ws.Range("G19").Value = 50.91514304475726
$ws.Range("H19").Value = 19.68788602811756
$ws.Range("I19").Value = 31.58893535041548
$ws.Range("J19").Value = 9.432036391997272
$ws.Range("L19").Value = 9.4081977795484
$ws.Range("M19").Value = 70.34455287010249

$ws.Range("C20").Value = 6.658348439901342
$ws.Range("D20").Value = 5.895467056976933
$ws.Range("E20").Value = 11.55162618046133
$ws.Range("F20").Value = 40.31176301613653
$ws.Range("G20").Value = 50.69288793204867
$ws.Range("H20").Value = 19.60961154673912
$ws.Range("I20").Value = 31.43896562705179
$ws.Range("J20").Value = 9.42136442218087
$ws.Range("L20").Value = 9.414751053998234
$ws.Range("M20").Value = 71.46337666475316

$ws.Range("C21").Value = 6.656163206379627
$ws.Range("D21").Value = 5.958549403720129
$ws.Range("E21").Value = 11.65595944660433
$ws.Range("F21").Value = 39.66176178387069
$ws.Range("G21").Value = 49.99727007875332
$ws.Range("H21").Value = 19.35682755157421
$ws.Range("I21").Value = 30.95539126904125
$ws.Range("J21").Value = 9.389157691203074
$ws.Range("L21").Value = 9.439800688557295
$ws.Range("M21").Value = 75.09624108407249

$ws.Range("C22").Value = 6.655646595888384
$ws.Range("D22").Value = 5.999344301645884
$ws.Range("E22").Value = 11.72660151542082
$ws.Range("F22").Value = 39.25905999096706
$ws.Range("G22").Value = 49.58223224620185
$ws.Range("H22").Value = 19.19935990727041
$ws.Range("I22").Value = 30.65483501554014
$ws.Range("J22").Value = 9.370882619294928
$ws.Range("L22").Value = 9.45845286981414
$ws.Range("M22").Value = 77.37596944642796

$ws.Range("C23").Value = 6.65583887400572
$ws.Range("D23").Value = 5.977612802533137
$ws.Range("E23").Value = 11.68867612538726
$ws.Range("F23").Value = 39.47192734383522
$ws.Range("G23").Value = 49.800025900674
$ws.Range("H23").Value = 19.28268612041381
$ws.Range("I23").Value = 30.81380668324053
$ws.Range("J23").Value = 9.380379032462431
$ws.Range("L23").Value = 9.448287788055048
$ws.Range("M23").Value = 76.16791175597658

$ws.Range("C24").Value = 6.658399537485544
$ws.Range("D24").Value = 5.894429121169073
$ws.Range("E24").Value = 11.54996281970891
$ws.Range("F24").Value = 40.32275182333478
$ws.Range("G24").Value = 50.70490445168261
$ws.Range("H24").Value = 19.61387193326795
$ws.Range("I24").Value = 31.44712573342002
$ws.Range("J24").Value = 9.421936879842049
$ws.Range("L24").Value = 9.414379829335532
$ws.Range("M24").Value = 71.40240959023322

$ws.Range("C25").Value = 6.665071832573463
$ws.Range("D25").Value = 5.803035765368956
$ws.Range("E25").Value = 11.41108691131745
$ws.Range("F25").Value = 41.33177093493654
$ws.Range("G25").Value = 51.84150765758749
$ws.Range("H25").Value = 20.00350507092281
$ws.Range("I25").Value = 32.19452390848097
$ws.Range("J25").Value = 9.478329958313356
$ws.Range("L25").Value = 9.387457817955481
$ws.Range("M25").Value = 65.85972914127518
